# This script re-orders the per-trial stimulus rows (rows 2-41, columns G..V)
# according to a fixed permutation, and renumbers the sequential trial_total
# column (F) starting at 406. Columns A-E (subject_id, task, block_total,
# block_scene, trial_block) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of new row number -> source (old) row number, for the block of
# data rows 2..41. This is the permutation observed between the original
# and target stimulus orderings.
$mapping = @{
    2  = 10
    3  = 30
    4  = 15
    5  = 39
    6  = 25
    7  = 3
    8  = 20
    9  = 14
    10 = 4
    11 = 13
    12 = 22
    13 = 9
    14 = 21
    15 = 7
    16 = 12
    17 = 17
    18 = 32
    19 = 31
    20 = 8
    21 = 40
    22 = 36
    23 = 41
    24 = 26
    25 = 19
    26 = 28
    27 = 24
    28 = 23
    29 = 33
    30 = 27
    31 = 5
    32 = 6
    33 = 34
    34 = 29
    35 = 11
    36 = 38
    37 = 18
    38 = 37
    39 = 35
    40 = 2
    41 = 16
}

$firstRow = 2
$lastRow = 41

# Snapshot the original contents of columns G..V (7..22) for every data row
# before any writes happen, since the permutation both reads from and
# writes to this same range.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 7; $c -le 22; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Apply the permutation: new row $r gets the G..V content that used to live
# in row $mapping[$r]. Also renumber column F (6) sequentially starting
# at 406.
$newTrialTotal = 406
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r]
    $srcVals = $snapshot[$srcRow]

    $ws.Cells.Item($r, 6).Value = $newTrialTotal
    $newTrialTotal++

    for ($c = 7; $c -le 22; $c++) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
